# Kalyankar -> Pusapati grading pass
# Fills in "Total Points" (earned) and grading comments for the
# "toString() method" rubric rows (rows 6 and 14), updates the running
# totals, and resets the sheet's scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Customer Class block (rows 3-6) ---------------------------------
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "(-1) for wrong output format"

# --- Product Class block (rows 10-14) --------------------------------
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "(-1) for wrong output format"

# --- Recalculate so the cached totals (E7, E15, E38, ...) are fresh --
$excel.Calculate()

# --- Reset view: scroll back to top, select F14 ----------------------
$ws.Activate()
$ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
